$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.916.17'
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").Value = '3.564.91'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").Value = '209.18'
$ws.Range("E5").Value = '  +11.43%  '

$ws.Range("D6").Value = '565.32'
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("D7").Value = '3.557.46'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '0.612'
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("D11").Value = '61.39'
$ws.Range("E11").Value = '  +10.45%  '

$ws.Range("D12").Value = '0.147'
$ws.Range("E12").Value = '  -1.94%  '

$ws.Range("E13").Value = '  +5.06%  '

$ws.Range("D14").Value = '10.18'

$ws.Range("D15").Value = '4.137.33'
$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("D16").Value = '3.561.66'
$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").Value = '18.89'
$ws.Range("E18").Value = '  +3.02%  '

$ws.Range("D19").Value = '67.691.87'
$ws.Range("E19").Value = '  +1.64%  '

$ws.Range("D20").Value = '12.12'
$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("D21").Value = '1.06'
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("D22").Value = '400.98'
$ws.Range("E22").Value = '  +1.45%  '

$ws.Range("D23").Value = '12.59'
$ws.Range("E23").Value = '  +11.66%  '

$ws.Range("D24").Value = '4.14'
$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("D25").Value = '83.89'
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("D26").Value = '2.87'
$ws.Range("E26").Value = '  -1.88%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '3.89'
$ws.Range("E27").Value = '  +8.25%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '12.35'
$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("D29").Value = '9.17'
$ws.Range("E29").Value = '  +3.37%  '

$ws.Range("D30").Value = '7.63'
$ws.Range("E30").Value = '  +1.06%  '

$ws.Range("D31").Value = '31.38'
$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("D32").Value = "'660.90"
$ws.Range("E32").Value = '  +3.63%  '

$ws.Range("D33").Value = '12.03'
$ws.Range("E33").Value = '  -1.19%  '

$ws.Range("D34").Value = '63.05'
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("E35").Value = '  -1.32%  '

$ws.Range("D36").Value = '40.94'
$ws.Range("E36").Value = '  -2.59%  '

$ws.Range("D37").Value = '0.408'
$ws.Range("E37").Value = '  +1.51%  '

$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("D39").Value = '3.24'
$ws.Range("E39").Value = '  +9.68%  '

$ws.Range("D40").Value = '0.0₃0747'
$ws.Range("E40").Value = '  -0.70%  '

$ws.Range("D41").Value = '3.155.89'
$ws.Range("E41").Value = '  +0.76%  '

$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("E44").Value = '  -0.11%  '

$ws.Range("E45").Value = '  +11.99%  '

$ws.Range("D46").Value = '0.0408'
$ws.Range("E46").Value = '  -1.05%  '

$ws.Range("D47").Value = "'0.130"
$ws.Range("E47").Value = '  -0.40%  '

$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '8.62'
$ws.Range("E48").Value = '  +2.16%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.57'
$ws.Range("E49").Value = '  +7.57%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = "'3.00"
$ws.Range("E50").Value = '  -2.58%  '

$ws.Range("D51").Value = '137.94'
$ws.Range("E51").Value = '  -1.30%  '
